$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) retains its text formatting so that
# numeric-looking values (e.g. "568.70", "11.30") are not silently
# coerced into numbers (which would drop formatting / trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '62.904.45'
$ws.Range('E2').Value = '  -1.47%  '
$ws.Range('D3').Value = '2.536.57'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').Value = '568.70'
$ws.Range('E5').Value = '  -0.98%  '
$ws.Range('D6').Value = '145.44'
$ws.Range('E6').Value = '  -0.73%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('D8').Value = '0.581'
$ws.Range('E8').Value = '  -1.24%  '
$ws.Range('D9').Value = '2.535.01'
$ws.Range('E9').Value = '  -0.07%  '
$ws.Range('E10').Value = '  -1.55%  '
$ws.Range('D11').Value = '5.48'
$ws.Range('E11').Value = '  -5.02%  '
$ws.Range('E12').Value = '  -0.46%  '
$ws.Range('D13').Value = '0.353'
$ws.Range('E13').Value = '  -1.82%  '
$ws.Range('D14').Value = '27.10'
$ws.Range('E14').Value = '  -2.84%  '
$ws.Range('D15').Value = '2.987.05'
$ws.Range('E15').Value = '  -0.27%  '
$ws.Range('D16').Value = '62.803.94'
$ws.Range('E16').Value = '  -1.30%  '
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').Value = '2.532.67'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('D19').Value = '11.26'
$ws.Range('E19').Value = '  -1.64%  '
$ws.Range('D20').Value = '333.43'
$ws.Range('E20').Value = '  -2.93%  '
$ws.Range('D21').Value = '4.31'
$ws.Range('E21').Value = '  -0.99%  '
$ws.Range('D22').Value = '6.72'
$ws.Range('E22').Value = '  -2.45%  '
$ws.Range('E23').Value = '  -0.64%  '
$ws.Range('D24').Value = '65.26'
$ws.Range('E24').Value = '  -0.95%  '
$ws.Range('D25').Value = '0.170'
$ws.Range('E25').Value = '  -1.46%  '
$ws.Range('D26').Value = '1.59'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('E27').Value = '  -0.14%  '
$ws.Range('D28').Value = '8.31'
$ws.Range('E28').Value = '  +0.85%  '
$ws.Range('D29').Value = '1.45'
$ws.Range('E29').Value = '  +1.10%  '
$ws.Range('D30').Value = '7.25'
$ws.Range('E30').Value = '  +6.44%  '
$ws.Range('D31').Value = '0.0₃0807'
$ws.Range('E32').Value = '  -2.11%  '
$ws.Range('D33').Value = '176.45'
$ws.Range('E33').Value = '  -0.22%  '
$ws.Range('E34').Value = '  -0.66%  '
$ws.Range('D35').Value = '400.42'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').Value = '19.05'
$ws.Range('D37').Value = '0.398'
$ws.Range('E37').Value = '  -1.08%  '
$ws.Range('E38').Value = '  +0.00%  '
$ws.Range('D39').Value = '4.31'
$ws.Range('E39').Value = '  -2.67%  '
$ws.Range('E40').Value = '  -0.18%  '
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('D42').Value = '39.45'
$ws.Range('E42').Value = '  -3.62%  '
$ws.Range('D43').Value = '150.90'
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').Value = '3.71'
$ws.Range('E44').Value = '  -2.04%  '
$ws.Range('D45').Value = '20.66'
$ws.Range('E45').Value = '  -1.88%  '
$ws.Range('D46').Value = '0.0529'
$ws.Range('E46').Value = '  -0.28%  '
$ws.Range('D47').Value = '0.597'
$ws.Range('E47').Value = '  -2.27%  '
$ws.Range('D48').Value = '0.0962'
$ws.Range('E48').Value = '  -0.51%  '
$ws.Range('D49').Value = '0.0237'
$ws.Range('E49').Value = '  +2.22%  '
$ws.Range('D50').Value = '18.11'
$ws.Range('E50').Value = '  -4.70%  '
$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D51').Value = '11.30'
$ws.Range('E51').Value = '  +0.34%  '

# Restore the original (default/Normal) style on the Price column now
# that the text values are safely stored, so no stray number-format
# style lingers on the cells.
$ws.Range("D2:D51").Style = "Normal"
